$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Cells.Item(2, 4) "26.921.66"
$ws.Cells.Item(2, 5).Value = "  +2.09%  "
Set-TextValue $ws.Cells.Item(3, 4) "1.650.38"
$ws.Cells.Item(3, 5).Value = "  +2.82%  "
$ws.Cells.Item(4, 5).Value = "  +0.09%  "
Set-TextValue $ws.Cells.Item(5, 4) "214.77"
$ws.Cells.Item(5, 5).Value = "  +1.30%  "
$ws.Cells.Item(6, 5).Value = "  +2.47%  "
$ws.Cells.Item(8, 5).Value = "  +2.61%  "
$ws.Cells.Item(9, 5).Value = "  +1.52%  "
Set-TextValue $ws.Cells.Item(10, 4) "20.16"
$ws.Cells.Item(10, 5).Value = "  +4.72%  "
$ws.Cells.Item(11, 5).Value = "  +2.30%  "
Set-TextValue $ws.Cells.Item(12, 4) "1.885.54"
$ws.Cells.Item(12, 5).Value = "  +2.98%  "
Set-TextValue $ws.Cells.Item(13, 4) "1.652.82"
$ws.Cells.Item(13, 5).Value = "  +2.76%  "
$ws.Cells.Item(14, 5).Value = "  +1.82%  "
$ws.Cells.Item(15, 5).Value = "  +2.45%  "
Set-TextValue $ws.Cells.Item(16, 4) "65.02"
$ws.Cells.Item(16, 5).Value = "  +2.74%  "
Set-TextValue $ws.Cells.Item(17, 4) "26.933.60"
$ws.Cells.Item(17, 5).Value = "  +2.15%  "
Set-TextValue $ws.Cells.Item(18, 4) "235.35"
$ws.Cells.Item(18, 5).Value = "  +2.31%  "
$ws.Cells.Item(19, 5).Value = "  +1.12%  "
Set-TextValue $ws.Cells.Item(20, 4) "7.71"
$ws.Cells.Item(20, 5).Value = "  +0.55%  "
$ws.Cells.Item(21, 5).Value = "  +0.02%  "
$ws.Cells.Item(22, 5).Value = "  +3.20%  "
$ws.Cells.Item(23, 5).Value = "  +3.86%  "
$ws.Cells.Item(24, 5).Value = "  +2.69%  "
Set-TextValue $ws.Cells.Item(25, 4) "145.30"
$ws.Cells.Item(25, 5).Value = "  -1.10%  "
$ws.Cells.Item(26, 5).Value = "  +1.96%  "
$ws.Cells.Item(27, 5).Value = "  +1.02%  "
$ws.Cells.Item(28, 5).Value = "  -0.01%  "
Set-TextValue $ws.Cells.Item(29, 4) "15.78"
$ws.Cells.Item(29, 5).Value = "  +2.40%  "
$ws.Cells.Item(30, 5).Value = "  +0.33%  "
$ws.Cells.Item(31, 5).Value = "  +1.63%  "
Set-TextValue $ws.Cells.Item(32, 4) "1.547.87"
$ws.Cells.Item(32, 5).Value = "  +3.92%  "
$ws.Cells.Item(33, 5).Value = "  +2.56%  "
$ws.Cells.Item(34, 5).Value = "  +4.64%  "
Set-TextValue $ws.Cells.Item(35, 4) "1.61"
$ws.Cells.Item(35, 5).Value = "  +8.95%  "
Set-TextValue $ws.Cells.Item(36, 4) "2.42"
$ws.Cells.Item(36, 5).Value = "  -0.19%  "
$ws.Cells.Item(37, 5).Value = "  +3.83%  "
Set-TextValue $ws.Cells.Item(38, 4) "0.890"
$ws.Cells.Item(38, 5).Value = "  +8.45%  "
$ws.Cells.Item(39, 5).Value = "  +2.81%  "
$ws.Cells.Item(40, 5).Value = "  +3.20%  "
$ws.Cells.Item(41, 5).Value = "  -0.01%  "
$ws.Cells.Item(42, 2).Value = "Aave"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Cells.Item(42, 4) "65.72"
$ws.Cells.Item(42, 5).Value = "  +8.02%  "
$ws.Cells.Item(43, 2).Value = "MXToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Cells.Item(43, 4) "2.24"
$ws.Cells.Item(43, 5).Value = "  +2.19%  "
Set-TextValue $ws.Cells.Item(44, 4) "1.792.58"
$ws.Cells.Item(44, 5).Value = "  +2.87%  "
Set-TextValue $ws.Cells.Item(45, 4) "0.775"
$ws.Cells.Item(45, 5).Value = "  +2.09%  "
Set-TextValue $ws.Cells.Item(46, 4) "0.926"
$ws.Cells.Item(46, 5).Value = "  -0.96%  "
Set-TextValue $ws.Cells.Item(47, 4) "90.06"
$ws.Cells.Item(47, 5).Value = "  +0.86%  "
$ws.Cells.Item(48, 5).Value = "  +1.60%  "
Set-TextValue $ws.Cells.Item(49, 4) "0.0989"
$ws.Cells.Item(50, 5).Value = "  +0.85%  "
Set-TextValue $ws.Cells.Item(51, 4) "7.61"
$ws.Cells.Item(51, 5).Value = "  +2.00%  "
